$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new End Time value for row 47 (the shared formula in D47 will recalc automatically)
$ws.Range("C47").Value = 0.086805555555555566

# Update the active selection to F43 (as recorded when the edit was made)
$ws.Range("F43").Select()
